$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common text values reused across rows
$title = "Tarea para departamento de sistema"
$descDigitel = "Se necesita configurar el ambiente de desarrollo para digitel"
$descFrontend = "Se necesita subir al repositorio los cambios del Front-end"
$users = "gabriel1407, user_admin"
$dept = "Departamento de prueba"

# Data rows 2-8: A, C(description), I(department), J(start_day), K(end_day)
$rows = @(
    @{ r = 2;  a = 35; c = $descDigitel;  j = 45145; k = 45146 },
    @{ r = 3;  a = 34; c = $descDigitel;  j = 45145; k = 45146 },
    @{ r = 4;  a = 33; c = $descDigitel;  j = 45145; k = 45146 },
    @{ r = 5;  a = 32; c = $descDigitel;  j = 45145; k = 45146 },
    @{ r = 6;  a = 31; c = $descDigitel;  j = 45145; k = 45146 },
    @{ r = 7;  a = 30; c = $descDigitel;  j = 45145; k = 45146 },
    @{ r = 8;  a = 29; c = $descFrontend; j = 45139; k = 45140 }
)

foreach ($row in $rows) {
    $r = $row.r

    $ws.Cells.Item($r, 1).Value = $row.a
    $ws.Cells.Item($r, 2).Value = $title
    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = $true
    $ws.Cells.Item($r, 7).Value = $true
    $ws.Cells.Item($r, 8).Value = $users
    $ws.Cells.Item($r, 9).Value = $dept

    $ws.Cells.Item($r, 10).Value = $row.j
    $ws.Cells.Item($r, 10).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($r, 11).Value = $row.k
    $ws.Cells.Item($r, 11).NumberFormat = "yyyy-mm-dd"
}

Write-Host "done"
